# Auto-generated edit script: applies numeric corrections to the
# H:N "profit" columns across several sheets (ALC, ARM, BSM, CRP,
# CUL, GSM, LTW, WVR) per the scheduled-runner data refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(135, 8).Value = 1355.7667  # H135
$ws.Cells.Item(135, 9).Value = 862.6  # I135
$ws.Cells.Item(135, 10).Value = 3821.6  # J135
$ws.Cells.Item(135, 11).Value = 7763.400000000001  # K135
$ws.Cells.Item(135, 12).Value = 34394.4  # L135
$ws.Cells.Item(135, 13).Value = -5228.400000000001  # M135
$ws.Cells.Item(135, 14).Value = -39464.4  # N135

$ws.Cells.Item(136, 8).Value = 34995  # H136
$ws.Cells.Item(136, 10).Value = 34995  # J136
$ws.Cells.Item(136, 12).Value = 34995  # L136
$ws.Cells.Item(136, 14).Value = -45195  # N136

$ws.Cells.Item(139, 8).Value = 0  # H139
$ws.Cells.Item(139, 10).Value = 0  # J139
$ws.Cells.Item(139, 12).Value = 0  # L139
$ws.Cells.Item(139, 14).ClearContents()  # N139

$ws.Cells.Item(140, 8).Value = 69573.63  # H140
$ws.Cells.Item(140, 10).Value = 67176  # J140
$ws.Cells.Item(140, 12).Value = 67176  # L140
$ws.Cells.Item(140, 14).Value = -77536  # N140

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 11581.321  # H32
$ws.Cells.Item(32, 9).Value = 12031.672  # I32
$ws.Cells.Item(32, 11).Value = 12031.672  # K32
$ws.Cells.Item(32, 13).Value = -11744.672  # M32

$ws.Cells.Item(133, 8).Value = 45003  # H133
$ws.Cells.Item(133, 10).Value = 45003  # J133
$ws.Cells.Item(133, 12).Value = 45003  # L133
$ws.Cells.Item(133, 14).Value = -50063  # N133

$ws.Cells.Item(138, 8).Value = 56092.25  # H138
$ws.Cells.Item(138, 10).Value = 56092.25  # J138
$ws.Cells.Item(138, 12).Value = 56092.25  # L138
$ws.Cells.Item(138, 14).Value = -66372.25  # N138

$ws.Cells.Item(139, 8).Value = 89392.42999999999  # H139
$ws.Cells.Item(139, 10).Value = 89392.42999999999  # J139
$ws.Cells.Item(139, 12).Value = 89392.42999999999  # L139
$ws.Cells.Item(139, 14).Value = -99672.42999999999  # N139

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(105, 8).Value = 3587.6667  # H105
$ws.Cells.Item(105, 9).Value = 3055.2  # I105
$ws.Cells.Item(105, 10).Value = 6250  # J105
$ws.Cells.Item(105, 11).Value = 3055.2  # K105
$ws.Cells.Item(105, 12).Value = 6250  # L105
$ws.Cells.Item(105, 13).Value = -1308.2  # M105
$ws.Cells.Item(105, 14).Value = -9744  # N105

$ws.Cells.Item(112, 8).Value = 21200  # H112
$ws.Cells.Item(112, 10).Value = 21200  # J112
$ws.Cells.Item(112, 12).Value = 21200  # L112
$ws.Cells.Item(112, 14).Value = -24154  # N112

$ws.Cells.Item(132, 8).Value = 75827.914  # H132
$ws.Cells.Item(132, 10).Value = 75827.914  # J132
$ws.Cells.Item(132, 12).Value = 75827.914  # L132
$ws.Cells.Item(132, 14).Value = -85947.914  # N132

$ws.Cells.Item(135, 8).Value = 56457.312  # H135
$ws.Cells.Item(135, 10).Value = 56457.312  # J135
$ws.Cells.Item(135, 12).Value = 56457.312  # L135
$ws.Cells.Item(135, 14).Value = -66597.31200000001  # N135

$ws.Cells.Item(137, 8).Value = 50000  # H137
$ws.Cells.Item(137, 10).Value = 50000  # J137
$ws.Cells.Item(137, 12).Value = 50000  # L137
$ws.Cells.Item(137, 14).Value = -60200  # N137

$ws.Cells.Item(138, 8).Value = 41379.6  # H138
$ws.Cells.Item(138, 10).Value = 41379.6  # J138
$ws.Cells.Item(138, 12).Value = 41379.6  # L138
$ws.Cells.Item(138, 14).Value = -51659.6  # N138

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 3128.3274  # H31
$ws.Cells.Item(31, 9).Value = 2330.625  # I31
$ws.Cells.Item(31, 10).Value = 4238.174  # J31
$ws.Cells.Item(31, 11).Value = 2330.625  # K31
$ws.Cells.Item(31, 12).Value = 4238.174  # L31
$ws.Cells.Item(31, 13).Value = -2035.625  # M31
$ws.Cells.Item(31, 14).Value = -4828.174  # N31

$ws.Cells.Item(34, 8).Value = 3128.3274  # H34
$ws.Cells.Item(34, 9).Value = 2330.625  # I34
$ws.Cells.Item(34, 10).Value = 4238.174  # J34
$ws.Cells.Item(34, 11).Value = 2330.625  # K34
$ws.Cells.Item(34, 12).Value = 4238.174  # L34
$ws.Cells.Item(34, 13).Value = -2128.625  # M34
$ws.Cells.Item(34, 14).Value = -4642.174  # N34

$ws.Cells.Item(122, 8).Value = 1169.35  # H122
$ws.Cells.Item(122, 9).Value = 1227.0714  # I122
$ws.Cells.Item(122, 10).Value = 1034.6666  # J122
$ws.Cells.Item(122, 11).Value = 3681.2142  # K122
$ws.Cells.Item(122, 12).Value = 3103.9998  # L122
$ws.Cells.Item(122, 13).Value = -1231.2142  # M122
$ws.Cells.Item(122, 14).Value = -8003.9998  # N122

$ws.Cells.Item(133, 8).Value = 0  # H133
$ws.Cells.Item(133, 10).Value = 0  # J133
$ws.Cells.Item(133, 12).Value = 0  # L133
$ws.Cells.Item(133, 14).ClearContents()  # N133

$ws.Cells.Item(135, 8).Value = 66316.42999999999  # H135
$ws.Cells.Item(135, 10).Value = 92047.78  # J135
$ws.Cells.Item(135, 12).Value = 92047.78  # L135
$ws.Cells.Item(135, 14).Value = -102187.78  # N135

$ws.Cells.Item(137, 8).Value = 52152.5  # H137
$ws.Cells.Item(137, 10).Value = 74305  # J137
$ws.Cells.Item(137, 12).Value = 74305  # L137
$ws.Cells.Item(137, 14).Value = -84505  # N137

$ws.Cells.Item(138, 8).Value = 52723.47  # H138
$ws.Cells.Item(138, 10).Value = 52723.47  # J138
$ws.Cells.Item(138, 12).Value = 52723.47  # L138
$ws.Cells.Item(138, 14).Value = -63003.47  # N138

$ws.Cells.Item(140, 8).Value = 71613.75  # H140
$ws.Cells.Item(140, 10).Value = 71613.75  # J140
$ws.Cells.Item(140, 12).Value = 71613.75  # L140
$ws.Cells.Item(140, 14).Value = -81973.75  # N140

$ws.Cells.Item(141, 8).Value = 80000  # H141
$ws.Cells.Item(141, 9).Value = 0  # I141
$ws.Cells.Item(141, 10).Value = 80000  # J141
$ws.Cells.Item(141, 11).Value = 0  # K141
$ws.Cells.Item(141, 12).Value = 80000  # L141
$ws.Cells.Item(141, 13).ClearContents()  # M141
$ws.Cells.Item(141, 14).Value = -90360  # N141

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(92, 8).Value = 1210.25  # H92
$ws.Cells.Item(92, 9).Value = 901  # I92
$ws.Cells.Item(92, 10).Value = 1313.3334  # J92
$ws.Cells.Item(92, 11).Value = 2703  # K92
$ws.Cells.Item(92, 12).Value = 3940.0002  # L92
$ws.Cells.Item(92, 13).Value = -1455  # M92
$ws.Cells.Item(92, 14).Value = -6436.0002  # N92

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(133, 8).Value = 60130  # H133
$ws.Cells.Item(133, 10).Value = 60130  # J133
$ws.Cells.Item(133, 12).Value = 60130  # L133
$ws.Cells.Item(133, 14).Value = -70250  # N133

$ws.Cells.Item(134, 8).Value = 35894  # H134
$ws.Cells.Item(134, 10).Value = 35894  # J134
$ws.Cells.Item(134, 12).Value = 107682  # L134
$ws.Cells.Item(134, 14).Value = -112752  # N134

$ws.Cells.Item(135, 8).Value = 48667  # H135
$ws.Cells.Item(135, 10).Value = 48667  # J135
$ws.Cells.Item(135, 12).Value = 48667  # L135
$ws.Cells.Item(135, 14).Value = -58807  # N135

$ws.Cells.Item(140, 8).Value = 49284.332  # H140
$ws.Cells.Item(140, 10).Value = 49284.332  # J140
$ws.Cells.Item(140, 12).Value = 49284.332  # L140
$ws.Cells.Item(140, 14).Value = -59644.332  # N140

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(122, 8).Value = 11369109  # H122
$ws.Cells.Item(122, 9).Value = 12505420  # I122
$ws.Cells.Item(122, 11).Value = 37516260  # K122
$ws.Cells.Item(122, 13).Value = -37513810  # M122

$ws.Cells.Item(132, 8).Value = 4671.952  # H132
$ws.Cells.Item(132, 9).Value = 4963.1035  # I132
$ws.Cells.Item(132, 10).Value = 4022.4614  # J132
$ws.Cells.Item(132, 11).Value = 14889.3105  # K132
$ws.Cells.Item(132, 12).Value = 12067.3842  # L132
$ws.Cells.Item(132, 13).Value = -12359.3105  # M132
$ws.Cells.Item(132, 14).Value = -17127.3842  # N132

$ws.Cells.Item(133, 8).Value = 74994.75  # H133
$ws.Cells.Item(133, 10).Value = 74994.75  # J133
$ws.Cells.Item(133, 12).Value = 74994.75  # L133
$ws.Cells.Item(133, 14).Value = -80054.75  # N133

$ws.Cells.Item(134, 8).Value = 75363.625  # H134
$ws.Cells.Item(134, 10).Value = 75363.625  # J134
$ws.Cells.Item(134, 12).Value = 75363.625  # L134
$ws.Cells.Item(134, 14).Value = -85503.625  # N134

$ws.Cells.Item(139, 8).Value = 37751.332  # H139
$ws.Cells.Item(139, 10).Value = 37751.332  # J139
$ws.Cells.Item(139, 12).Value = 37751.332  # L139
$ws.Cells.Item(139, 14).Value = -48031.332  # N139

$ws.Cells.Item(140, 8).Value = 70535.75  # H140
$ws.Cells.Item(140, 10).Value = 70535.75  # J140
$ws.Cells.Item(140, 12).Value = 70535.75  # L140
$ws.Cells.Item(140, 14).Value = -80895.75  # N140

$ws.Cells.Item(141, 8).Value = 45410.555  # H141
$ws.Cells.Item(141, 10).Value = 45410.555  # J141
$ws.Cells.Item(141, 12).Value = 45410.555  # L141
$ws.Cells.Item(141, 14).Value = -55770.555  # N141

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(46, 8).Value = 48352.668  # H46
$ws.Cells.Item(46, 10).Value = 48352.668  # J46
$ws.Cells.Item(46, 12).Value = 48352.668  # L46
$ws.Cells.Item(46, 14).Value = -48814.668  # N46

$ws.Cells.Item(133, 8).Value = 42194.2  # H133
$ws.Cells.Item(133, 10).Value = 42194.2  # J133
$ws.Cells.Item(133, 12).Value = 42194.2  # L133
$ws.Cells.Item(133, 14).Value = -52314.2  # N133

$ws.Cells.Item(134, 8).Value = 48352.668  # H134
$ws.Cells.Item(134, 10).Value = 48352.668  # J134
$ws.Cells.Item(134, 12).Value = 145058.004  # L134
$ws.Cells.Item(134, 14).Value = -150128.004  # N134

$ws.Cells.Item(135, 8).Value = 58191.363  # H135
$ws.Cells.Item(135, 10).Value = 58191.363  # J135
$ws.Cells.Item(135, 12).Value = 58191.363  # L135
$ws.Cells.Item(135, 14).Value = -68331.363  # N135

$ws.Cells.Item(137, 8).Value = 35747.25  # H137
$ws.Cells.Item(137, 10).Value = 35747.25  # J137
$ws.Cells.Item(137, 12).Value = 35747.25  # L137
$ws.Cells.Item(137, 14).Value = -45947.25  # N137

$ws.Cells.Item(139, 8).Value = 57843.332  # H139
$ws.Cells.Item(139, 10).Value = 57843.332  # J139
$ws.Cells.Item(139, 12).Value = 57843.332  # L139
$ws.Cells.Item(139, 14).Value = -68123.33199999999  # N139

$ws.Cells.Item(141, 8).Value = 81475.56  # H141
$ws.Cells.Item(141, 10).Value = 81475.56  # J141
$ws.Cells.Item(141, 12).Value = 81475.56  # L141
$ws.Cells.Item(141, 14).Value = -91835.56  # N141
